# ---------------------------------------------------------------------------
# Roll the balance sheet forward by one fiscal year (1396/12 drops off the
# front, 1401/12 is appended at the end) and refresh the "published on"
# dates to match. All numeric series in columns D:H shift one column to the
# left (D<-E, E<-F, F<-G, G<-H) and column H receives the newly reported
# 1401/12 figures.
# ---------------------------------------------------------------------------
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# -- Header: fiscal-period labels (row 8) and report publish dates (row 9) --
$ws.Cells.Item(8, 4).Value = "12 ماهه منتهی به 1397/12"
$ws.Cells.Item(8, 5).Value = "12 ماهه منتهی به 1398/12"
$ws.Cells.Item(8, 6).Value = "12 ماهه منتهی به 1399/12"
$ws.Cells.Item(8, 7).Value = "12 ماهه منتهی به 1400/12"
$ws.Cells.Item(8, 8).Value = "12 ماهه منتهی به 1401/12"
$ws.Cells.Item(9, 4).Value = "1399-04-21 (9)"
$ws.Cells.Item(9, 5).Value = "1400-04-27 (13)"
$ws.Cells.Item(9, 6).Value = "1401-04-18 (10)"
$ws.Cells.Item(9, 7).Value = "1402-02-27 (7)"
$ws.Cells.Item(9, 8).Value = "1402-02-27"

# -- Balance sheet figures: shift D:H left one column, append new 1401/12 --
$ws.Cells.Item(12, 4).Value = 691288
$ws.Cells.Item(12, 5).Value = 1630106
$ws.Cells.Item(12, 6).Value = 795824
$ws.Cells.Item(12, 7).Value = 1037258
$ws.Cells.Item(12, 8).Value = 1664766
$ws.Cells.Item(13, 4).Value = 684705
$ws.Cells.Item(13, 5).Value = 1081705
$ws.Cells.Item(13, 6).Value = 534948
$ws.Cells.Item(13, 7).Value = 1159948
$ws.Cells.Item(13, 8).Value = 3254948
$ws.Cells.Item(14, 4).Value = 1818477
$ws.Cells.Item(14, 5).Value = 2867698
$ws.Cells.Item(14, 6).Value = 5876340
$ws.Cells.Item(14, 7).Value = 13490026
$ws.Cells.Item(14, 8).Value = 27149784
$ws.Cells.Item(15, 4).Value = 13695561
$ws.Cells.Item(15, 5).Value = 16523076
$ws.Cells.Item(15, 6).Value = 29547150
$ws.Cells.Item(15, 7).Value = 40401383
$ws.Cells.Item(15, 8).Value = 51753043
$ws.Cells.Item(16, 4).Value = 352897
$ws.Cells.Item(16, 5).Value = 328528
$ws.Cells.Item(16, 6).Value = 311893
$ws.Cells.Item(16, 7).Value = 279257
$ws.Cells.Item(16, 8).Value = 2767720
$ws.Cells.Item(17, 4).Value = 0
$ws.Cells.Item(17, 5).Value = 0
$ws.Cells.Item(17, 6).Value = 0
$ws.Cells.Item(17, 7).Value = 0
$ws.Cells.Item(17, 8).Value = 0
$ws.Cells.Item(18, 4).Value = 17242928
$ws.Cells.Item(18, 5).Value = 22431113
$ws.Cells.Item(18, 6).Value = 37066155
$ws.Cells.Item(18, 7).Value = 56367872
$ws.Cells.Item(18, 8).Value = 86590261
$ws.Cells.Item(19, 4).Value = 18507
$ws.Cells.Item(19, 5).Value = 18347
$ws.Cells.Item(19, 6).Value = 12779
$ws.Cells.Item(19, 7).Value = 6818
$ws.Cells.Item(19, 8).Value = 9713
$ws.Cells.Item(20, 4).Value = 1000
$ws.Cells.Item(20, 5).Value = 1000
$ws.Cells.Item(20, 6).Value = 1000
$ws.Cells.Item(20, 7).Value = 1000
$ws.Cells.Item(20, 8).Value = 1000
$ws.Cells.Item(21, 4).Value = 0
$ws.Cells.Item(21, 5).Value = 0
$ws.Cells.Item(21, 6).Value = 0
$ws.Cells.Item(21, 7).Value = 0
$ws.Cells.Item(21, 8).Value = 0
$ws.Cells.Item(22, 4).Value = 10233893
$ws.Cells.Item(22, 5).Value = 10373942
$ws.Cells.Item(22, 6).Value = 10603803
$ws.Cells.Item(22, 7).Value = 12867330
$ws.Cells.Item(22, 8).Value = 14852692
$ws.Cells.Item(23, 4).Value = 0
$ws.Cells.Item(23, 5).Value = 0
$ws.Cells.Item(23, 6).Value = 0
$ws.Cells.Item(23, 7).Value = 0
$ws.Cells.Item(23, 8).Value = 0
$ws.Cells.Item(24, 4).Value = "-"
$ws.Cells.Item(24, 5).Value = "-"
$ws.Cells.Item(24, 6).Value = "-"
$ws.Cells.Item(24, 7).Value = "-"
$ws.Cells.Item(24, 8).Value = "-"
$ws.Cells.Item(25, 4).Value = 301421
$ws.Cells.Item(25, 5).Value = 248549
$ws.Cells.Item(25, 6).Value = 259393
$ws.Cells.Item(25, 7).Value = 366464
$ws.Cells.Item(25, 8).Value = 516995
$ws.Cells.Item(26, 4).Value = 10554821
$ws.Cells.Item(26, 5).Value = 10641838
$ws.Cells.Item(26, 6).Value = 10876975
$ws.Cells.Item(26, 7).Value = 13241612
$ws.Cells.Item(26, 8).Value = 15380400
$ws.Cells.Item(27, 4).Value = 27797749
$ws.Cells.Item(27, 5).Value = 33072951
$ws.Cells.Item(27, 6).Value = 47943130
$ws.Cells.Item(27, 7).Value = 69609484
$ws.Cells.Item(27, 8).Value = 101970661
$ws.Cells.Item(29, 4).Value = 19791294
$ws.Cells.Item(29, 5).Value = 24562795
$ws.Cells.Item(29, 6).Value = 16383446
$ws.Cells.Item(29, 7).Value = 24276712
$ws.Cells.Item(29, 8).Value = 24709646
$ws.Cells.Item(30, 4).Value = "-"
$ws.Cells.Item(30, 5).Value = "-"
$ws.Cells.Item(30, 6).Value = "-"
$ws.Cells.Item(30, 7).Value = "-"
$ws.Cells.Item(30, 8).Value = "-"
$ws.Cells.Item(31, 4).Value = 142991
$ws.Cells.Item(31, 5).Value = 0
$ws.Cells.Item(31, 6).Value = 0
$ws.Cells.Item(31, 7).Value = 0
$ws.Cells.Item(31, 8).Value = 0
$ws.Cells.Item(32, 4).Value = 463288
$ws.Cells.Item(32, 5).Value = 236180
$ws.Cells.Item(32, 6).Value = 3009488
$ws.Cells.Item(32, 7).Value = 2111508
$ws.Cells.Item(32, 8).Value = 4638750
$ws.Cells.Item(33, 4).Value = 33118
$ws.Cells.Item(33, 5).Value = 448474
$ws.Cells.Item(33, 6).Value = 97544
$ws.Cells.Item(33, 7).Value = 165877
$ws.Cells.Item(33, 8).Value = 123591
$ws.Cells.Item(34, 4).Value = 2296475
$ws.Cells.Item(34, 5).Value = 1909768
$ws.Cells.Item(34, 6).Value = 4471208
$ws.Cells.Item(34, 7).Value = 4643795
$ws.Cells.Item(34, 8).Value = 8883658
$ws.Cells.Item(35, 4).Value = 0
$ws.Cells.Item(35, 5).Value = 0
$ws.Cells.Item(35, 6).Value = 0
$ws.Cells.Item(35, 7).Value = 1
$ws.Cells.Item(35, 8).Value = 0
$ws.Cells.Item(36, 4).Value = 0
$ws.Cells.Item(36, 5).Value = 0
$ws.Cells.Item(36, 6).Value = 0
$ws.Cells.Item(36, 7).Value = 0
$ws.Cells.Item(36, 8).Value = 0
$ws.Cells.Item(37, 4).Value = 22727166
$ws.Cells.Item(37, 5).Value = 27157217
$ws.Cells.Item(37, 6).Value = 23961686
$ws.Cells.Item(37, 7).Value = 31197893
$ws.Cells.Item(37, 8).Value = 38355645
$ws.Cells.Item(38, 4).Value = 0
$ws.Cells.Item(38, 5).Value = 0
$ws.Cells.Item(38, 6).Value = 0
$ws.Cells.Item(38, 7).Value = 0
$ws.Cells.Item(38, 8).Value = 0
$ws.Cells.Item(39, 4).Value = "-"
$ws.Cells.Item(39, 5).Value = "-"
$ws.Cells.Item(39, 6).Value = "-"
$ws.Cells.Item(39, 7).Value = "-"
$ws.Cells.Item(39, 8).Value = "-"
$ws.Cells.Item(40, 4).Value = 0
$ws.Cells.Item(40, 5).Value = 0
$ws.Cells.Item(40, 6).Value = 0
$ws.Cells.Item(40, 7).Value = 0
$ws.Cells.Item(40, 8).Value = 0
$ws.Cells.Item(41, 4).Value = 431001
$ws.Cells.Item(41, 5).Value = 543605
$ws.Cells.Item(41, 6).Value = 704866
$ws.Cells.Item(41, 7).Value = 957536
$ws.Cells.Item(41, 8).Value = 1198830
$ws.Cells.Item(42, 4).Value = 431001
$ws.Cells.Item(42, 5).Value = 543605
$ws.Cells.Item(42, 6).Value = 704866
$ws.Cells.Item(42, 7).Value = 957536
$ws.Cells.Item(42, 8).Value = 1198830
$ws.Cells.Item(43, 4).Value = 23158167
$ws.Cells.Item(43, 5).Value = 27700822
$ws.Cells.Item(43, 6).Value = 24666552
$ws.Cells.Item(43, 7).Value = 32155429
$ws.Cells.Item(43, 8).Value = 39554475
$ws.Cells.Item(45, 4).Value = 1143422
$ws.Cells.Item(45, 5).Value = 1143422
$ws.Cells.Item(45, 6).Value = 1143422
$ws.Cells.Item(45, 7).Value = 1143422
$ws.Cells.Item(45, 8).Value = 8000000
$ws.Cells.Item(46, 4).Value = 0
$ws.Cells.Item(46, 5).Value = 0
$ws.Cells.Item(46, 6).Value = 0
$ws.Cells.Item(46, 7).Value = 48667
$ws.Cells.Item(46, 8).Value = 101282
$ws.Cells.Item(47, 4).Value = 0
$ws.Cells.Item(47, 5).Value = 0
$ws.Cells.Item(47, 6).Value = 0
$ws.Cells.Item(47, 7).Value = 0
$ws.Cells.Item(47, 8).Value = 0
$ws.Cells.Item(48, 4).Value = 0
$ws.Cells.Item(48, 5).Value = 0
$ws.Cells.Item(48, 6).Value = -128406
$ws.Cells.Item(48, 7).Value = -146820
$ws.Cells.Item(48, 8).Value = -60407
$ws.Cells.Item(49, 4).Value = 0
$ws.Cells.Item(49, 5).Value = 0
$ws.Cells.Item(49, 6).Value = 829
$ws.Cells.Item(49, 7).Value = 0
$ws.Cells.Item(49, 8).Value = 0
$ws.Cells.Item(50, 4).Value = 114342
$ws.Cells.Item(50, 5).Value = 114342
$ws.Cells.Item(50, 6).Value = 114342
$ws.Cells.Item(50, 7).Value = 114342
$ws.Cells.Item(50, 8).Value = 114342
$ws.Cells.Item(51, 4).Value = 0
$ws.Cells.Item(51, 5).Value = 0
$ws.Cells.Item(51, 6).Value = 0
$ws.Cells.Item(51, 7).Value = 0
$ws.Cells.Item(51, 8).Value = 0
$ws.Cells.Item(52, 4).Value = "-"
$ws.Cells.Item(52, 5).Value = "-"
$ws.Cells.Item(52, 6).Value = "-"
$ws.Cells.Item(52, 7).Value = "-"
$ws.Cells.Item(52, 8).Value = "-"
$ws.Cells.Item(53, 4).Value = 0
$ws.Cells.Item(53, 5).Value = 0
$ws.Cells.Item(53, 6).Value = 0
$ws.Cells.Item(53, 7).Value = 0
$ws.Cells.Item(53, 8).Value = 0
$ws.Cells.Item(54, 4).Value = "-"
$ws.Cells.Item(54, 5).Value = "-"
$ws.Cells.Item(54, 6).Value = "-"
$ws.Cells.Item(54, 7).Value = "-"
$ws.Cells.Item(54, 8).Value = "-"
$ws.Cells.Item(55, 4).Value = 0
$ws.Cells.Item(55, 5).Value = 0
$ws.Cells.Item(55, 6).Value = 0
$ws.Cells.Item(55, 7).Value = 0
$ws.Cells.Item(55, 8).Value = 0
$ws.Cells.Item(56, 4).Value = 3381818
$ws.Cells.Item(56, 5).Value = 4114365
$ws.Cells.Item(56, 6).Value = 22146391
$ws.Cells.Item(56, 7).Value = 36294444
$ws.Cells.Item(56, 8).Value = 54260969
$ws.Cells.Item(57, 4).Value = 4639582
$ws.Cells.Item(57, 5).Value = 5372129
$ws.Cells.Item(57, 6).Value = 23276578
$ws.Cells.Item(57, 7).Value = 37454055
$ws.Cells.Item(57, 8).Value = 62416186
$ws.Cells.Item(58, 4).Value = 27797749
$ws.Cells.Item(58, 5).Value = 33072951
$ws.Cells.Item(58, 6).Value = 47943130
$ws.Cells.Item(58, 7).Value = 69609484
$ws.Cells.Item(58, 8).Value = 101970661
